# Helper: write a value as TEXT (matching the source inlineStr cells) while
# keeping the cell's style at the default "Normal" (no leftover custom style).
function Set-TextValue($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Helper: write a real number value.
function Set-NumberValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Helper: write a boolean value.
function Set-BoolValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "quotations": insert a new row at position 17 (shifts old rows
# 17-23 down to 18-24) and fill it with the new quotation record.
# ---------------------------------------------------------------------------
$wsQuot = $wb.Worksheets.Item("quotations")
$wsQuot.Rows("17:17").Insert()

Set-TextValue   $wsQuot 17 1  "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"
Set-TextValue   $wsQuot 17 2  "BKMVCD56XX"
Set-TextValue   $wsQuot 17 3  "MARCOS ROBINSON MARTINS DE SOUZA"
Set-TextValue   $wsQuot 17 4  ""
Set-TextValue   $wsQuot 17 5  ""
Set-BoolValue   $wsQuot 17 6  $false
Set-TextValue   $wsQuot 17 7  "40.894999999999996"
Set-TextValue   $wsQuot 17 8  "40.894999999999996"
Set-TextValue   $wsQuot 17 9  "Pendente"
Set-TextValue   $wsQuot 17 10 "2025-10-20T12:54:09.917Z"
Set-TextValue   $wsQuot 17 11 ""
Set-TextValue   $wsQuot 17 12 ""
Set-TextValue   $wsQuot 17 13 "Kaue Teixeira Caldeira Venâncio"
Set-TextValue   $wsQuot 17 14 ""
Set-TextValue   $wsQuot 17 15 "2025-10-13T13:10:02.576Z"
Set-TextValue   $wsQuot 17 16 ""
Set-TextValue   $wsQuot 17 17 "percentage"
Set-TextValue   $wsQuot 17 18 "0"
Set-TextValue   $wsQuot 17 19 "0"
Set-TextValue   $wsQuot 17 20 "NDgyNTIyNjo1NzAxNg=="
Set-TextValue   $wsQuot 17 21 "pending"

# ---------------------------------------------------------------------------
# Sheet "items": insert 5 new rows starting at position 57 (shifts old rows
# 57-85 down to 62-90) and fill them with the new line items belonging to
# the new quotation created above.
# ---------------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("items")
$wsItems.Rows("57:61").Insert()

# Row 57
Set-TextValue   $wsItems 57 1  "NGM4MjM2YzAtZWJkZC00NTY1LWFkNGEtNjBiYzEwYzVlMjM5OjU3MDE2"
Set-NumberValue $wsItems 57 2  2
Set-NumberValue $wsItems 57 3  96
Set-TextValue   $wsItems 57 4  ""
Set-NumberValue $wsItems 57 5  6
Set-TextValue   $wsItems 57 6  "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"
Set-TextValue   $wsItems 57 7  "ODcwZTI1ZDEtMTRkNC00M2IyLTk0MTItOGJhNDdiYzIzMjg1OjU3MDE2"
Set-NumberValue $wsItems 57 8  48
Set-TextValue   $wsItems 57 9  "product"
Set-TextValue   $wsItems 57 10 "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"

# Row 58
Set-TextValue   $wsItems 58 1  "NmJjNjE2Y2UtOGFjZi00OTdiLTgyMWEtYTI4NmMzMDA1MjViOjU3MDE2"
Set-NumberValue $wsItems 58 2  1
Set-NumberValue $wsItems 58 3  6000000000000001
Set-TextValue   $wsItems 58 4  ""
Set-NumberValue $wsItems 58 5  6
Set-TextValue   $wsItems 58 6  "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"
Set-TextValue   $wsItems 58 7  "OTgxYjZlMTAtNGZiMy00YjAwLWI4OTYtMTcxNGM5MTg2Y2NiOjU3MDE2"
Set-NumberValue $wsItems 58 8  6000000000000001
Set-TextValue   $wsItems 58 9  "product"
Set-TextValue   $wsItems 58 10 "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"

# Row 59
Set-TextValue   $wsItems 59 1  "YTYwMTIwMzItNGMzNy00MGYyLThiMjktY2FhNjAzZTljY2YyOjU3MDE2"
Set-NumberValue $wsItems 59 2  1
Set-NumberValue $wsItems 59 3  38825000000000000
Set-TextValue   $wsItems 59 4  ""
Set-NumberValue $wsItems 59 5  6
Set-TextValue   $wsItems 59 6  "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"
Set-TextValue   $wsItems 59 7  "YTQ4MzMzNDUtNWU2Yy00YmVmLWE4OWYtY2Y4ZWNkNjI2ZjM5OjU3MDE2"
Set-NumberValue $wsItems 59 8  38825000000000000
Set-TextValue   $wsItems 59 9  "product"
Set-TextValue   $wsItems 59 10 "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"

# Row 60
Set-TextValue   $wsItems 60 1  "YTk2NmJhMzgtNDJkYS00OTQ3LWJlZWYtNzVlOGQzOTA4NjM4OjU3MDE2"
Set-NumberValue $wsItems 60 2  1
Set-NumberValue $wsItems 60 3  21000000000000000
Set-TextValue   $wsItems 60 4  ""
Set-NumberValue $wsItems 60 5  6
Set-TextValue   $wsItems 60 6  "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"
Set-TextValue   $wsItems 60 7  "YmU4MzA1YWYtMThkMC00ODhkLTk0NTMtMzE3MGVkYjI2NGFiOjU3MDE2"
Set-NumberValue $wsItems 60 8  21000000000000000
Set-TextValue   $wsItems 60 9  "product"
Set-TextValue   $wsItems 60 10 "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"

# Row 61
Set-TextValue   $wsItems 61 1  "ZTI4OGEyZDItNTQ3Ni00MmVjLWFlZTMtZDk4N2U3OGU3ODJmOjU3MDE2"
Set-NumberValue $wsItems 61 2  1
Set-NumberValue $wsItems 61 3  30000000000000000
Set-TextValue   $wsItems 61 4  ""
Set-NumberValue $wsItems 61 5  6
Set-TextValue   $wsItems 61 6  "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"
Set-TextValue   $wsItems 61 7  "MWI1OWUzZTYtMTQxYi00NjY4LWFjZjAtNGI1OGE2Mjk3Nzk4OjU3MDE2"
Set-NumberValue $wsItems 61 8  30000000000000000
Set-TextValue   $wsItems 61 9  "product"
Set-TextValue   $wsItems 61 10 "OTc4YTljYTItZGQ1Ni00MDIyLWEyMTItYWQwZWY4M2VkNWYxOjU3MDE2"
